$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "false start" rows (rows 2 and 3). The rows that follow
# (old rows 4 and 5) shift up to become the new rows 2 and 3.
$ws.Rows("2:3").Select()
$ws.Rows("2:3").Delete()
